# edit.ps1 -- apply the tracked-changes-style edit described by the commit:
# "Changed Thread.sleep() to use .join() instead. Also added more test files
#  and uploaded new pdf version of report."
#
# This touches only the body paragraphs of the lab report: a series of small
# wording tweaks (several of which replace the Thread.sleep()-based
# explanation with a .join()-based one), plus justifying (both-align) most
# body paragraphs, plus a small rPr addition on the image run.

$d = $word.ActiveDocument

function Replace-Text {
    param(
        [string]$Find,
        [string]$ReplaceWith
    )
    $ok = $d.Content.Find.Execute(
        $Find, $true, $false, $false, $false, $false,
        $true, 1, $false, $ReplaceWith, 2)
    if (-not $ok) {
        Write-Output ("NOT FOUND: " + $Find)
    }
}

# 1. "Therefore, " removed / capitalization fix
Replace-Text "three classes. Therefore, an analysis" "three classes. An analysis"

# 2-4. FindDefective -> FindDefective(); sub-halves -> sub-arrays; "bulbs are
#      their" -> "bulbs at their". These all live inside the single run that
#      follows the spellStart/spellEnd-wrapped "FindDefective" word, so doing
#      them as one Replace-Text call leaves that pairing untouched.
Replace-Text " is used. Each call to this method splits the array into two sub-halves and creates a new thread for each half to further recurse through the sub-arrays. This recursion is repeated until all the defective lightbulbs are found. When the multi-threaded recursive search terminates, the program outputs the total number of defective bulbs are their respective original indices in the array" "() is used. Each call to this method splits the array into two sub-arrays and creates a new thread for each half to further recurse through the sub-arrays. This recursion is repeated until all the defective lightbulbs are found. When the multi-threaded recursive search terminates, the program outputs the total number of defective bulbs at their respective original indices in the array"

# 5. "missing; it is provided" -> "missing; if it is provided"
Replace-Text "argument is missing; it is provided" "argument is missing; if it is provided"

# 6. "make a call to " -> "make a call to the " (stop right before the
#    spellStart/gramStart-wrapped "GetBulbArray" run)
Replace-Text "make a call to " "make a call to the "

# 7-8. "information is used in the input .txt file. The " ->
#      "information is provided in the input .txt file. Finally, the "
#      (single run, stops right before the spellStart/gramStart-wrapped
#      "GetBulbArray" run)
Replace-Text "no erroneous information is used in the input .txt file. The " "no erroneous information is provided in the input .txt file. Finally, the "

# 9. "method finally returns" -> "method returns"
Replace-Text ") method finally returns the integer array" ") method returns the integer array"

# 10. Thread.sleep() explanation -> .join() explanation
Replace-Text " In the meantime, the main thread is set to sleep for a specified period of time until the threads find all of the defective bulbs in the array. Keep in mind" " Additionally, the main thread makes a call to the .join() method on the DefectiveThread object, which will wait until the sub-threads return before continuing in the main thread. Keep in mind"

# 11. "in order to break out from the recursion." -> "to break out from the recursive calls."
Replace-Text "is of size 1 in order to break out from the recursion. If this is not the case" "is of size 1 to break out from the recursive calls. If this is not the case"

# 12. Append new sentence about .join() at end of the recursive-method paragraph
# (match only the tail run "(), etc.)." so the FindDefective spellStart/spellEnd
#  pairing a few runs earlier is left untouched)
$rsquo = [char]0x2019
Replace-Text "(), etc.)." ("(), etc.). After the .start() method is called, a .join() method follows to ensure that the program pauses correctly and doesn" + $rsquo + "t bubble back up to the main program prematurely.")

# 13. "In order to ensure that there are no race conditions" -> "To ensure that..."
#     Match exactly the gramStart/gramEnd-wrapped "In order to" run (it is the
#     unique, capitalized occurrence -- the paragraph-opening one) so the
#     proofErr pairing around it is replaced cleanly instead of orphaned.
Replace-Text "In order to" "To"

# 14. "introduced by introducing static objects." -> "introduced by via static objects."
Replace-Text "synchronization is introduced by introducing static objects." "synchronization is introduced by via static objects."

# 15. "incrementation of defective lightbulb count" -> "incrementation of the defective lightbulb count"
Replace-Text "lie in the incrementation of defective lightbulb count" "lie in the incrementation of the defective lightbulb count"

# 16. "synchronize the aforementioned variables" -> "synchronize the variables"
Replace-Text "synchronize the aforementioned variables so that" "synchronize the variables so that"

# 17. "passing those objects to the static synchronized" -> "passing those objects to the synchronized"
Replace-Text "passing those objects to the static " "passing those objects to the "

# 18. "lock object in order to ensure" -> "lock object to ensure"
Replace-Text "lock object in order to ensure that" "lock object to ensure that"

# 19. "this operation that the incrementation" -> "this operation than the incrementation"
Replace-Text "this operation that the incrementation" "this operation than the incrementation"

# 20. "continue their operation once the main thread wakes from its sleep. As a final step, the"
#     -> "continue their remaining operations. At this point, the"
Replace-Text "continue their operation once the main thread wakes from its sleep. As a final step, the" "continue their remaining operations. At this point, the"

# 21. "class calls its " -> "class calls its static " (text immediately before
#     the GetBulbArray run, which is wrapped in spellStart/gramStart markers
#     that should stay intact)
Replace-Text "FileReader.java class calls its " "FileReader.java class calls its static "

# 22. final paragraph rewording about .join() calls / synchronized variables
Replace-Text "by creating new threads in order to search for the defective bulbs. When all of the defective bulbs are found," "by creating new threads to search for the defective bulbs. Along the way, a series of .join() methods are called on the new threads in order to allow for correct program flow. To ensure correct functionality of the program, several defective bulb counting variables are synchronized in order to bypass the race condition. When all the defective bulbs are found,"

# --- Paragraph justification: set Alignment = wdAlignParagraphJustify (3) on
#     the body paragraphs that gained <w:jc w:val="both"/> ---
$justifyStarts = @(
    "The program is comprised of three classes.",
    "First and foremost, this program is a command line program",
    "The main method and entry point of the program is through Program.java",
    "Next, back in Program.java, a",
    "To ensure that there are no race conditions",
    "The above process continues until",
    "When each thread completes as per the above specifications,",
    "The process is now summarized in a concluding statement."
)

foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text
    foreach ($prefix in $justifyStarts) {
        if ($ptext.StartsWith($prefix)) {
            $p.Range.ParagraphFormat.Alignment = 3
            break
        }
    }
}

# --- Add <w:rPr><w:noProof/></w:rPr> on the run holding the inline drawing ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.InlineShapes.Count -gt 0) {
        $p.Range.Font.NoProof = $true
    }
}

Write-Output "done"
